# Generate Report for Handoff
# The f0820fe7-175a-4264-a194-9d9298b01209.md row (row 3 on every sheet) moves
# from "Handed back: in sync with en-US" to "Ready for handoff", picks up a
# fresh handoff timestamp, and gets a new "error detail" note about the
# handback file being stale.

$wb = $excel.ActiveWorkbook

$statusText   = "Ready for handoff"
$errorDetail  = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/366d51502d629243ecd436162288f1cded126d23/e2e/f0820fe7-175a-4264-a194-9d9298b01209.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8fa29608c4613731990cabc11230164e4a2d02b6/e2e/f0820fe7-175a-4264-a194-9d9298b01209.md."

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = $statusText
$ws.Range("F3").Value = $statusText
$ws.Range("G3").Value = "2016-08-30 16:57:23"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = $statusText
$ws.Range("H3").Value = "2016-08-30 16:57:18"
$ws.Range("P3").Value = $errorDetail
# Stored OOXML <col> width = ColumnWidth + 5/6 for this engine, so set
# ColumnWidth a bit under 40 to land on a stored width of exactly 40.
$ws.Columns.Item(16).ColumnWidth = 39.16666666666667

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = $statusText
$ws.Range("H3").Value = "2016-08-30 16:57:23"
$ws.Range("P3").Value = $errorDetail
$ws.Columns.Item(16).ColumnWidth = 39.16666666666667
